$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "0191448213678"
$ws.Range("B4").Value = 1
